# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 07:52"

# --- India (row 22) ---
$ws.Range("B22").Value = 11555
$ws.Range("C22").Value = 68
$ws.Range("D22").Value = 1362
$ws.Range("E22").Value = 9797
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 396

# --- Barein (row 64) ---
$ws.Range("D64").Value = 648
$ws.Range("E64").Value = 873

# --- Lituania (row 73) ---
$ws.Range("D73").Value = 138
$ws.Range("E73").Value = 924

# --- Oman overtakes Republica de Macedonia / Camerun / Eslovaquia in the
#     case-count ranking, so rows 78-81 shift down one position each and
#     Oman's refreshed figures land in row 78 (Cuba in row 82 is untouched).
$ws.Range("A78").Value = "Oman"
$ws.Range("B78").Value = 910
$ws.Range("C78").Value = 97
$ws.Range("D78").Value = 130
$ws.Range("E78").Value = 776
$ws.Range("F78").Value = 3
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 4

$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("B79").Value = 908
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 86
$ws.Range("E79").Value = 778
$ws.Range("F79").Value = 15
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 44

$ws.Range("A80").Value = "Camerun"
$ws.Range("B80").Value = 848
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 130
$ws.Range("E80").Value = 704
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 14

$ws.Range("A81").Value = "Eslovaquia"
$ws.Range("B81").Value = 835
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 113
$ws.Range("E81").Value = 720
$ws.Range("F81").Value = 5
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 2

# --- Maldivas (row 172) ---
$ws.Range("B172").Value = 21
$ws.Range("C172").Value = 1
$ws.Range("E172").Value = 5
